$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update cell values (shared-string edits from the diff) ---
$ws.Range("B10").Value = "8426375 - Wendell de Queiróz Lamas"
$ws.Range("C10").Value = "8426375 - Wendell de Queiróz Lamas"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Energy sources and the Brazilian energy matrix. Anthropogenic activities, energy demand, and socio-economic development. Availability of sources and evaluation of energy generation potential. Electrical power fundamentals. Power plants. Solar energy. Wind energy. Fossil energy. Biomass energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development."
$ws.Range("C14").Value = "Energy sources and the Brazilian energy matrix. Anthropogenic activities, energy demand, and socio-economic development. Availability of sources and evaluation of energy generation potential. Electrical power fundamentals. Power plants. Solar energy. Wind energy. Fossil energy. Biomass energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development."
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2015"
$ws.Range("C15").Value = "01/01/2015"
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Energy sources and the Brazilian energy matrix. Anthropogenic activities, energy demand, and socio-economic development. Availability of sources and evaluation of energy generation potential. Electrical power fundamentals: generation, transmission, and distribution. Power plants: hydraulic, thermal, and nuclear. Solar energy. Wind energy. Fossil energy. Biomass energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development."
$ws.Range("C16").Value = "Energy sources and the Brazilian energy matrix. Anthropogenic activities, energy demand, and socio-economic development. Availability of sources and evaluation of energy generation potential. Electrical power fundamentals: generation, transmission, and distribution. Power plants: hydraulic, thermal, and nuclear. Solar energy. Wind energy. Fossil energy. Biomass energy. Environmental impacts of energy generation, transmission, availability, and supply in regional development."
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8426375 - Wendell de Queiróz Lamas"
$ws.Range("C18").Value = "8426375 - Wendell de Queiróz Lamas"
$ws.Range("A19").Value = "Critério:"
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B23").Value = "LOB1021 -  Física IV  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOB1021 -  Física IV  (Requisito fraco)`n"

# --- Clear cells that no longer hold content ---
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Range("A23").Clear()
$ws.Range("B24").Clear()
$ws.Range("C24").Clear()

# --- Delete now-empty trailing row 24 (dimension shrinks to C23) ---
$ws.Rows.Item(24).Delete()

# --- Row height adjustments ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
